# [Update] Joined both BOM files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "bom" -> "bom (1)"
$ws.Name = "bom (1)"

# A2: collapse the "RI1-3, RL1-108m " range shorthand into the fully
# enumerated designator list (both BOMs' RL/RI designators joined together).
$ws.Range("A2").Value = "RL108,RL107,RL106,RL105,RL104,RL103,RL102,RL101,RL100,RL99,RL98,RL97,RL96,RL95,RL94,RL93,RL92,RL91,RL90,RL89,RL88,RL87,RL86,RL85,RL84,RL83,RL82,RL81,RL80,RL79,RL78,RL77,RL76,RL75,RL74,RL73,RL72,RL71,RL70,RL69,RL68,RL67,RL66,RL65,RL64,RL63,RL62,RL60,RL58,RL57,RL56,RL55,RL54,RL53,RL52,RL51,RL50,RL49,RL48,RL47,RL46,RL45,RL44,RL43,RL42,RL41,RL40,RL39,RL38,RL37,RL36,RL34,RL33,RL32,RL31,RL30,RL29,RL28,RL27,RL26,RL25,RL24,RL23,RL22,RL21,RL20,RL19,RL18,RL17,RL16,RL15,RL14,RL13,RL12,RL11,RL10,RL9,RL8,RL7,RL6,RL5,RL4,RL3,RL2,RL1,RI3,RI2,RI1"

# A10: collapse the "D1-108, DRST1" range shorthand into the fully
# enumerated diode designator list.
$ws.Range("A10").Value = "D1,D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51,D61,D62,D63,D64,D65,D66,D67,D68,D69,D70,D71,D73,D79,D96,D80,D97,D81,D98,D82,D83,D84,D99,D85,D86,D87,D88,D100,D89,D101,D102,D13,D32,D52,D91,D103,D14,D33,D53,D15,D35,D54,D74,D104,D16,D36,D55,D75,D92,D105,D17,D37,D56,D76,D93,D106,D18,D38,D57,D77,D94,D107,D19,D39,D59,D95,D108,D72,D90,DRST1"

# E13: fuse's LCSC part number corrected C70117 -> C70118.
$ws.Range("E13").Value = "C70118"

# Page setup: switch the sheet to landscape printing with standard margins.
$ws.PageSetup.Orientation = 2
$ws.PageSetup.TopMargin = 54.0
$ws.PageSetup.BottomMargin = 54.0
$ws.PageSetup.LeftMargin = 50.4
$ws.PageSetup.RightMargin = 50.4
$ws.PageSetup.HeaderMargin = 0.0
$ws.PageSetup.FooterMargin = 0.0
